$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1702.5385
$ws.Range("I9").Value = 122.166664
$ws.Range("J9").Value = 3057.1428
$ws.Range("K9").Value = 122.166664
$ws.Range("L9").Value = 3057.1428
$ws.Range("M9").Value = 46.833336
$ws.Range("N9").Value = -3395.1428
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -20968
$ws.Range("H58").Value = 975
$ws.Range("H62").Value = 5697.2
$ws.Range("I62").Value = 6503.5
$ws.Range("K62").Value = 6503.5
$ws.Range("M62").Value = -5879.5
$ws.Range("H65").Value = 5697.2
$ws.Range("I65").Value = 6503.5
$ws.Range("K65").Value = 32517.5
$ws.Range("M65").Value = -29397.5
$ws.Range("H80").Value = 2300.25
$ws.Range("I80").Value = 1680
$ws.Range("J80").Value = 2743.2856
$ws.Range("K80").Value = 5040
$ws.Range("L80").Value = 8229.856800000001
$ws.Range("M80").Value = -4042
$ws.Range("N80").Value = -10225.8568
$ws.Range("H81").Value = 78124.875
$ws.Range("J81").Value = 78124.875
$ws.Range("L81").Value = 78124.875
$ws.Range("N81").Value = -80120.875
$ws.Range("H83").Value = 2300.25
$ws.Range("I83").Value = 1680
$ws.Range("J83").Value = 2743.2856
$ws.Range("K83").Value = 15120
$ws.Range("L83").Value = 24689.5704
$ws.Range("M83").Value = -10128
$ws.Range("N83").Value = -34673.5704
$ws.Range("H84").Value = 78124.875
$ws.Range("J84").Value = 78124.875
$ws.Range("L84").Value = 234374.625
$ws.Range("N84").Value = -244358.625
$ws.Range("H87").Value = 95000
$ws.Range("J87").Value = 95000
$ws.Range("L87").Value = 95000
$ws.Range("N87").Value = -97496
$ws.Range("H90").Value = 95000
$ws.Range("J90").Value = 95000
$ws.Range("L90").Value = 285000
$ws.Range("N90").Value = -297480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 13499.5
$ws.Range("I61").Value = 19999
$ws.Range("K61").Value = 19999
$ws.Range("M61").Value = -19787
$ws.Range("H101").Value = 1799.5
$ws.Range("J101").Value = 1799.5
$ws.Range("L101").Value = 1799.5
$ws.Range("N101").Value = -8289.5
$ws.Range("H132").Value = 632
$ws.Range("I132").Value = 647.6667
$ws.Range("K132").Value = 1943.0001
$ws.Range("M132").Value = 586.9999
$ws.Range("H136").Value = 13499.5
$ws.Range("I136").Value = 19999
$ws.Range("K136").Value = 59997
$ws.Range("M136").Value = -57447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5608.364
$ws.Range("I20").Value = 2086.875
$ws.Range("J20").Value = 14999
$ws.Range("K20").Value = 2086.875
$ws.Range("L20").Value = 14999
$ws.Range("M20").Value = -1839.875
$ws.Range("N20").Value = -15493
$ws.Range("H86").Value = 2521.7778
$ws.Range("I86").Value = 2529.4285
$ws.Range("K86").Value = 2529.4285
$ws.Range("M86").Value = -1406.4285
$ws.Range("H89").Value = 2521.7778
$ws.Range("I89").Value = 2529.4285
$ws.Range("K89").Value = 12647.1425
$ws.Range("M89").Value = -7031.1425
$ws.Range("H94").Value = 1798.6666
$ws.Range("I94").Value = 1158.4
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1158.4
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -707.4000000000001
$ws.Range("N94").Value = -5902
$ws.Range("H134").Value = 1346.6428
$ws.Range("I134").Value = 1334.8462
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 4004.5386
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -1469.5386
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2997.4546
$ws.Range("I31").Value = 3049
$ws.Range("J31").Value = 2765.5
$ws.Range("K31").Value = 3049
$ws.Range("L31").Value = 2765.5
$ws.Range("M31").Value = -2754
$ws.Range("N31").Value = -3355.5
$ws.Range("H34").Value = 2997.4546
$ws.Range("I34").Value = 3049
$ws.Range("J34").Value = 2765.5
$ws.Range("K34").Value = 3049
$ws.Range("L34").Value = 2765.5
$ws.Range("M34").Value = -2847
$ws.Range("N34").Value = -3169.5
$ws.Range("H86").Value = 7747885
$ws.Range("I86").Value = 9960138
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 9960138
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -9959015
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 7747885
$ws.Range("I89").Value = 9960138
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 49800690
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -49795074
$ws.Range("N89").Value = -36232
$ws.Range("H105").Value = 1383.3334
$ws.Range("I105").Value = 1202.5
$ws.Range("J105").Value = 1745
$ws.Range("K105").Value = 1202.5
$ws.Range("L105").Value = 1745
$ws.Range("M105").Value = 544.5
$ws.Range("N105").Value = -5239

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 110123
$ws.Range("I2").Value = 156.25
$ws.Range("J2").Value = 183434.17
$ws.Range("K2").Value = 937.5
$ws.Range("L2").Value = 1100605.02
$ws.Range("M2").Value = -824.5
$ws.Range("N2").Value = -1100831.02
$ws.Range("H23").Value = 991.25
$ws.Range("I23").Value = 980
$ws.Range("K23").Value = 2940
$ws.Range("M23").Value = -2705
$ws.Range("H41").Value = 1083.4286
$ws.Range("I41").Value = 1255.6666
$ws.Range("J41").Value = 954.25
$ws.Range("K41").Value = 3766.9998
$ws.Range("L41").Value = 2862.75
$ws.Range("M41").Value = -3428.9998
$ws.Range("N41").Value = -3538.75
$ws.Range("H80").Value = 2500
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2500
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H87").Value = 616.3333
$ws.Range("I87").Value = 616.3333
$ws.Range("K87").Value = 1848.9999
$ws.Range("M87").Value = -600.9999
$ws.Range("H90").Value = 616.3333
$ws.Range("I90").Value = 616.3333
$ws.Range("K90").Value = 5546.9997
$ws.Range("M90").Value = 693.0002999999997
$ws.Range("H114").Value = 1250
$ws.Range("I114").Value = 1250
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 3750
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -496
$ws.Range("N114").ClearContents()
$ws.Range("H117").Value = 450
$ws.Range("I117").Value = 450
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1350
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2092
$ws.Range("N117").ClearContents()
$ws.Range("H134").Value = 4599.75
$ws.Range("I134").Value = 4699.727
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 14099.181
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -9029.181
$ws.Range("N134").Value = -20640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 959.1667
$ws.Range("I2").Value = 1130.1111
$ws.Range("J2").Value = 446.33334
$ws.Range("K2").Value = 1130.1111
$ws.Range("L2").Value = 446.33334
$ws.Range("M2").Value = -1017.1111
$ws.Range("N2").Value = -672.33334
$ws.Range("H97").Value = 1688.5555
$ws.Range("I97").Value = 1399.875
$ws.Range("K97").Value = 1399.875
$ws.Range("M97").Value = -903.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 55000
$ws.Range("J76").Value = 55000
$ws.Range("L76").Value = 55000
$ws.Range("N76").Value = -55630
$ws.Range("H79").Value = 55000
$ws.Range("J79").Value = 55000
$ws.Range("L79").Value = 55000
$ws.Range("N79").Value = -57184
